$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix label text (remove stray trailing spaces)
$ws.Range("A1").Value = "email"

# Update the email value and turn it into a mailto hyperlink
$ws.Range("B1").Value = "abc155553@test.com"
$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:abc155553@test.com", "", "", "abc155553@test.com")

$ws.Range("A11").Value = "City"

# Column width adjustments
$ws.Columns.Item(1).ColumnWidth = 14.28515625
$ws.Columns.Item(2).ColumnWidth = 15

# Selection moves to A12
$ws.Range("A12").Select()
